# Add role add messages
# - Rename "get_content" message to "get_content_list" (rows 27 & 28)
# - Fill in parameter names for the get_content_list call (row 27 & 28)
# - Add a new "approve_content" message pair (rows 29 & 30)
# - Add a new "add_role_define" message pair (rows 31 & 32)
# - Add a new "add_role" message pair (rows 33 & 34)
# - Move selection/scroll to reflect the last-edited area (row 35)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename existing message name "get_content" -> "get_content_list"
# (both request & response rows reference the same message name)
$ws.Range("B27").Value = "get_content_list"
$ws.Range("B28").Value = "get_content_list"

# Request parameters for get_content_list (row 27)
$ws.Range("E27").Value = "Type ID"

# Response parameters for get_content_list (row 28)
$ws.Range("E28").Value = "Message status"
$ws.Range("F28").Value = "Content list"

# New message: approve_content - request (row 29)
$ws.Range("B29").Value = "approve_content"
$ws.Range("C29").Value = "Web Client"
$ws.Range("D29").Value = "Server"
$ws.Range("E29").Value = "Session Code"
$ws.Range("F29").Value = "Content ID"

# New message: approve_content - response (row 30)
$ws.Range("B30").Value = "approve_content"
$ws.Range("C30").Value = "Server"
$ws.Range("D30").Value = "Web Client"
$ws.Range("E30").Value = "Message status"

# New message: add_role_define - request (row 31)
$ws.Range("B31").Value = "add_role_define"
$ws.Range("C31").Value = "Web Client"
$ws.Range("D31").Value = "Server"
$ws.Range("E31").Value = "Session Code"
$ws.Range("F31").Value = "Role name"

# New message: add_role_define - response (row 32)
$ws.Range("B32").Value = "add_role_define"
$ws.Range("C32").Value = "Server"
$ws.Range("D32").Value = "Web Client"
$ws.Range("E32").Value = "Message status"

# New message: add_role - request (row 33)
$ws.Range("B33").Value = "add_role"
$ws.Range("C33").Value = "Web Client"
$ws.Range("D33").Value = "Server"
$ws.Range("E33").Value = "Session Code"
$ws.Range("F33").Value = "Role ID"
$ws.Range("G33").Value = "Type ID"
$ws.Range("H33").Value = "Role Code"

# New message: add_role - response (row 34)
$ws.Range("B34").Value = "add_role"
$ws.Range("C34").Value = "Server"
$ws.Range("D34").Value = "Web Client"
$ws.Range("E34").Value = "Message status"

# Reflect the final cursor/scroll position left by the edit session
$window = $excel.ActiveWindow
$window.ScrollRow = 10
$window.ScrollColumn = 1
$ws.Range("E35").Select()
